# Auto-applies numeric corrections to the Leve profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1213.12
$ws.Range("J19").Value = 2693
$ws.Range("L19").Value = 2693
$ws.Range("N19").Value = -3043

# Row 100
$ws.Range("H100").Value = 1367.1538
$ws.Range("I100").Value = 1434.6364
$ws.Range("J100").Value = 996
$ws.Range("K100").Value = 1434.6364
$ws.Range("L100").Value = 996
$ws.Range("M100").Value = -893.6364000000001
$ws.Range("N100").Value = -2078

# Row 129
$ws.Range("H129").Value = 1911.1
$ws.Range("J129").Value = 2102.0667
$ws.Range("L129").Value = 6306.2001
$ws.Range("N129").Value = -16306.2001

# Row 132
$ws.Range("H132").Value = 2242.3333
$ws.Range("I132").Value = 1286.875
$ws.Range("J132").Value = 9886
$ws.Range("K132").Value = 3860.625
$ws.Range("L132").Value = 29658
$ws.Range("M132").Value = -1330.625
$ws.Range("N132").Value = -34718

# Row 137
$ws.Range("H137").Value = 1248.1111
$ws.Range("I137").Value = 779.125
$ws.Range("K137").Value = 2337.375
$ws.Range("M137").Value = 212.625

# Row 138
$ws.Range("H138").Value = 3053.628
$ws.Range("I138").Value = 2038.5555
$ws.Range("J138").Value = 3322.3235
$ws.Range("K138").Value = 6115.666499999999
$ws.Range("L138").Value = 9966.970499999999
$ws.Range("M138").Value = -975.6664999999994
$ws.Range("N138").Value = -20246.9705

# Row 141
$ws.Range("H141").Value = 6749.5
$ws.Range("I141").Value = 3999.3333
$ws.Range("K141").Value = 11997.9999
$ws.Range("M141").Value = -6817.999899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 1849
$ws.Range("I3").Value = 1849
$ws.Range("K3").Value = 1849
$ws.Range("M3").Value = -1734

# Row 122
$ws.Range("H122").Value = 8001.316
$ws.Range("I122").Value = 5819.2583
$ws.Range("J122").Value = 17664.715
$ws.Range("K122").Value = 17457.7749
$ws.Range("L122").Value = 52994.145
$ws.Range("M122").Value = -15007.7749
$ws.Range("N122").Value = -57894.145

# Row 123
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3575.0908
$ws.Range("I20").Value = 2432.6
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 2432.6
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -2185.6
$ws.Range("N20").Value = -15494

# Row 33
$ws.Range("H33").Value = 12250
$ws.Range("I33").Value = 9000
$ws.Range("K33").Value = 9000
$ws.Range("M33").Value = -8664

# Row 105
$ws.Range("H105").Value = 3627103.2
$ws.Range("I105").Value = 5558944.5
$ws.Range("J105").Value = 4900.875
$ws.Range("K105").Value = 5558944.5
$ws.Range("L105").Value = 4900.875
$ws.Range("M105").Value = -5557197.5
$ws.Range("N105").Value = -8394.875

# Row 134
$ws.Range("H134").Value = 3024.5
$ws.Range("I134").Value = 2403.5833
$ws.Range("K134").Value = 7210.749899999999
$ws.Range("M134").Value = -4675.749899999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 64780.81
$ws.Range("I22").Value = 86742.64
$ws.Range("J22").Value = 20857.143
$ws.Range("K22").Value = 86742.64
$ws.Range("L22").Value = 20857.143
$ws.Range("M22").Value = -86392.64
$ws.Range("N22").Value = -21557.143

# Row 59
$ws.Range("H59").Value = 61999.332
$ws.Range("J59").Value = 61999.332
$ws.Range("L59").Value = 61999.332
$ws.Range("N59").Value = -64289.332

# Row 141
$ws.Range("H141").Value = 120000
$ws.Range("J141").Value = 120000
$ws.Range("L141").Value = 120000
$ws.Range("N141").Value = -130360

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 3292469.8
$ws.Range("I4").Value = 4354399
$ws.Range("J4").Value = 488.9
$ws.Range("K4").Value = 13063197
$ws.Range("L4").Value = 1466.7
$ws.Range("M4").Value = -13063085
$ws.Range("N4").Value = -1690.7

# Row 46
$ws.Range("H46").Value = 2500701
$ws.Range("J46").Value = 3333935
$ws.Range("L46").Value = 10001805
$ws.Range("N46").Value = -10001987

$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 7199.8
$ws.Range("I36").Value = 7999.6665
$ws.Range("J36").Value = 6000
$ws.Range("K36").Value = 7999.6665
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -7514.6665
$ws.Range("N36").Value = -6970

# Row 107
$ws.Range("H107").Value = 1340.2222
$ws.Range("I107").Value = 236.6
$ws.Range("K107").Value = 236.6
$ws.Range("M107").Value = 1683.4

# Row 132
$ws.Range("H132").Value = 4698.5557
$ws.Range("I132").Value = 4698.5557
$ws.Range("K132").Value = 14095.6671
$ws.Range("M132").Value = -11565.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2203
$ws.Range("I7").Value = 2376.6428
$ws.Range("J7").Value = 987.5
$ws.Range("K7").Value = 2376.6428
$ws.Range("L7").Value = 987.5
$ws.Range("M7").Value = -2264.6428
$ws.Range("N7").Value = -1211.5

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 122
$ws.Range("H122").Value = 6985.3335
$ws.Range("I122").Value = 6945
$ws.Range("K122").Value = 20835
$ws.Range("M122").Value = -18385

# Row 126
$ws.Range("H126").Value = 2203
$ws.Range("I126").Value = 2376.6428
$ws.Range("J126").Value = 987.5
$ws.Range("K126").Value = 7129.928400000001
$ws.Range("L126").Value = 2962.5
$ws.Range("M126").Value = -4659.928400000001
$ws.Range("N126").Value = -7902.5

# Row 132
$ws.Range("H132").Value = 202400.8
$ws.Range("I132").Value = 202400.8
$ws.Range("K132").Value = 607202.3999999999
$ws.Range("M132").Value = -604672.3999999999

# Row 136
$ws.Range("H136").Value = 5443.8667
$ws.Range("I136").Value = 4750.1
$ws.Range("J136").Value = 6831.4
$ws.Range("K136").Value = 14250.3
$ws.Range("L136").Value = 20494.2
$ws.Range("M136").Value = -11700.3
$ws.Range("N136").Value = -25594.2

$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 17995
$ws.Range("I12").Value = 17995
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 17995
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -17853
$ws.Range("N12").ClearContents()

# Row 45
$ws.Range("H45").Value = 20712.5
$ws.Range("I45").Value = 17600
$ws.Range("K45").Value = 17600
$ws.Range("M45").Value = -17109

# Row 96
$ws.Range("H96").Value = 1933
$ws.Range("I96").Value = 1933
$ws.Range("K96").Value = 1933
$ws.Range("M96").Value = -560

# Row 122
$ws.Range("H122").Value = 3136.087
$ws.Range("I122").Value = 2934.1428
$ws.Range("J122").Value = 3450.2222
$ws.Range("K122").Value = 8802.428400000001
$ws.Range("L122").Value = 10350.6666
$ws.Range("M122").Value = -6352.428400000001
$ws.Range("N122").Value = -15250.6666
